$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header X2: "Utility" -> "Utility (Percent)"
$ws.Range("X2").Value = "Utility (Percent)"

# Append " usec" to latency values in columns L, M, N, O, P, Q for rows 3-38
foreach ($row in 3..38) {
    foreach ($col in @("L", "M", "N", "O", "P", "Q")) {
        $cell = $ws.Range("$col$row")
        $cell.Value = "$($cell.Value2) usec"
    }
}
